$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of |S*|/n column (J), bold size 11
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Re-use B14's formatting for B15:B17 (copy/paste-format keeps the
# stylesheet from accumulating transient/orphan cellXfs entries that
# sequential per-property font mutation on each cell would create).
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Selection as saved in the final workbook
$ws.Range("J12").Select()

# Print setup seen in the final file
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
